$wb = $excel.ActiveWorkbook

# Sheet 1: LP1912 -- update header + rows 42-75
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 07:19:37"
$ws1.Range("A3").Value = "Total filas: 70"
$ws1.Cells.Item(42, 1).Value = "07:19:37"
$ws1.Cells.Item(42, 2).Value = "07:19"
$ws1.Cells.Item(42, 3).Value = "10_OLMOS"
$ws1.Cells.Item(42, 4).Value = 0
$ws1.Cells.Item(42, 5).Value = "LP1912"
$ws1.Cells.Item(43, 1).Value = "05:49:10"
$ws1.Cells.Item(43, 2).Value = "07:21"
$ws1.Cells.Item(43, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(43, 4).Value = 92
$ws1.Cells.Item(43, 5).Value = "LP1912"
$ws1.Cells.Item(44, 1).Value = "06:14:19"
$ws1.Cells.Item(44, 2).Value = "07:23"
$ws1.Cells.Item(44, 3).Value = "10_OLMOS"
$ws1.Cells.Item(44, 4).Value = 69
$ws1.Cells.Item(44, 5).Value = "LP1912"
$ws1.Cells.Item(45, 1).Value = "06:14:19"
$ws1.Cells.Item(45, 2).Value = "07:31"
$ws1.Cells.Item(45, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(45, 4).Value = 77
$ws1.Cells.Item(45, 5).Value = "LP1912"
$ws1.Cells.Item(46, 1).Value = "07:19:37"
$ws1.Cells.Item(46, 2).Value = "07:31"
$ws1.Cells.Item(46, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(46, 4).Value = 12
$ws1.Cells.Item(46, 5).Value = "LP1912"
$ws1.Cells.Item(47, 1).Value = "05:49:10"
$ws1.Cells.Item(47, 2).Value = "07:32"
$ws1.Cells.Item(47, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(47, 4).Value = 103
$ws1.Cells.Item(47, 5).Value = "LP1912"
$ws1.Cells.Item(48, 1).Value = "05:49:10"
$ws1.Cells.Item(48, 2).Value = "07:32"
$ws1.Cells.Item(48, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(48, 4).Value = 103
$ws1.Cells.Item(48, 5).Value = "LP1912"
$ws1.Cells.Item(49, 1).Value = "05:49:10"
$ws1.Cells.Item(49, 2).Value = "07:32"
$ws1.Cells.Item(49, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(49, 4).Value = 103
$ws1.Cells.Item(49, 5).Value = "LP1912"
$ws1.Cells.Item(50, 1).Value = "07:19:37"
$ws1.Cells.Item(50, 2).Value = "07:34"
$ws1.Cells.Item(50, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(50, 4).Value = 15
$ws1.Cells.Item(50, 5).Value = "LP1912"
$ws1.Cells.Item(51, 1).Value = "07:19:37"
$ws1.Cells.Item(51, 2).Value = "07:36"
$ws1.Cells.Item(51, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(51, 4).Value = 17
$ws1.Cells.Item(51, 5).Value = "LP1912"
$ws1.Cells.Item(52, 1).Value = "05:49:10"
$ws1.Cells.Item(52, 2).Value = "07:37"
$ws1.Cells.Item(52, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(52, 4).Value = 108
$ws1.Cells.Item(52, 5).Value = "LP1912"
$ws1.Cells.Item(53, 1).Value = "05:49:10"
$ws1.Cells.Item(53, 2).Value = "07:39"
$ws1.Cells.Item(53, 3).Value = "10_OLMOS"
$ws1.Cells.Item(53, 4).Value = 110
$ws1.Cells.Item(53, 5).Value = "LP1912"
$ws1.Cells.Item(54, 1).Value = "06:14:19"
$ws1.Cells.Item(54, 2).Value = "07:47"
$ws1.Cells.Item(54, 3).Value = "14_ABASTO"
$ws1.Cells.Item(54, 4).Value = 93
$ws1.Cells.Item(54, 5).Value = "LP1912"
$ws1.Cells.Item(55, 1).Value = "05:49:10"
$ws1.Cells.Item(55, 2).Value = "07:48"
$ws1.Cells.Item(55, 3).Value = "14_ABASTO"
$ws1.Cells.Item(55, 4).Value = 119
$ws1.Cells.Item(55, 5).Value = "LP1912"
$ws1.Cells.Item(56, 1).Value = "06:14:19"
$ws1.Cells.Item(56, 2).Value = "07:51"
$ws1.Cells.Item(56, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(56, 4).Value = 97
$ws1.Cells.Item(56, 5).Value = "LP1912"
$ws1.Cells.Item(57, 1).Value = "06:43:12"
$ws1.Cells.Item(57, 2).Value = "07:52"
$ws1.Cells.Item(57, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(57, 4).Value = 69
$ws1.Cells.Item(57, 5).Value = "LP1912"
$ws1.Cells.Item(58, 1).Value = "07:19:37"
$ws1.Cells.Item(58, 2).Value = "07:59"
$ws1.Cells.Item(58, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(58, 4).Value = 40
$ws1.Cells.Item(58, 5).Value = "LP1912"
$ws1.Cells.Item(59, 1).Value = "06:14:19"
$ws1.Cells.Item(59, 2).Value = "08:00"
$ws1.Cells.Item(59, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(59, 4).Value = 106
$ws1.Cells.Item(59, 5).Value = "LP1912"
$ws1.Cells.Item(60, 1).Value = "07:19:37"
$ws1.Cells.Item(60, 2).Value = "08:03"
$ws1.Cells.Item(60, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(60, 4).Value = 44
$ws1.Cells.Item(60, 5).Value = "LP1912"
$ws1.Cells.Item(61, 1).Value = "06:43:12"
$ws1.Cells.Item(61, 2).Value = "08:03"
$ws1.Cells.Item(61, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(61, 4).Value = 80
$ws1.Cells.Item(61, 5).Value = "LP1912"
$ws1.Cells.Item(62, 1).Value = "06:57:11"
$ws1.Cells.Item(62, 2).Value = "08:05"
$ws1.Cells.Item(62, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(62, 4).Value = 68
$ws1.Cells.Item(62, 5).Value = "LP1912"
$ws1.Cells.Item(63, 1).Value = "07:19:37"
$ws1.Cells.Item(63, 2).Value = "08:10"
$ws1.Cells.Item(63, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(63, 4).Value = 51
$ws1.Cells.Item(63, 5).Value = "LP1912"
$ws1.Cells.Item(64, 1).Value = "06:14:19"
$ws1.Cells.Item(64, 2).Value = "08:12"
$ws1.Cells.Item(64, 3).Value = "15_ABASTO"
$ws1.Cells.Item(64, 4).Value = 118
$ws1.Cells.Item(64, 5).Value = "LP1912"
$ws1.Cells.Item(65, 1).Value = "06:43:12"
$ws1.Cells.Item(65, 2).Value = "08:21"
$ws1.Cells.Item(65, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(65, 4).Value = 98
$ws1.Cells.Item(65, 5).Value = "LP1912"
$ws1.Cells.Item(66, 1).Value = "07:19:37"
$ws1.Cells.Item(66, 2).Value = "08:22"
$ws1.Cells.Item(66, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(66, 4).Value = 63
$ws1.Cells.Item(66, 5).Value = "LP1912"
$ws1.Cells.Item(67, 1).Value = "06:43:12"
$ws1.Cells.Item(67, 2).Value = "08:23"
$ws1.Cells.Item(67, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(67, 4).Value = 100
$ws1.Cells.Item(67, 5).Value = "LP1912"
$ws1.Cells.Item(68, 1).Value = "06:43:12"
$ws1.Cells.Item(68, 2).Value = "08:23"
$ws1.Cells.Item(68, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(68, 4).Value = 100
$ws1.Cells.Item(68, 5).Value = "LP1912"
$ws1.Cells.Item(69, 1).Value = "06:43:12"
$ws1.Cells.Item(69, 2).Value = "08:27"
$ws1.Cells.Item(69, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(69, 4).Value = 104
$ws1.Cells.Item(69, 5).Value = "LP1912"
$ws1.Cells.Item(70, 1).Value = "06:57:11"
$ws1.Cells.Item(70, 2).Value = "08:42"
$ws1.Cells.Item(70, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(70, 4).Value = 105
$ws1.Cells.Item(70, 5).Value = "LP1912"
$ws1.Cells.Item(71, 1).Value = "07:19:37"
$ws1.Cells.Item(71, 2).Value = "08:43"
$ws1.Cells.Item(71, 3).Value = "14_ABASTO"
$ws1.Cells.Item(71, 4).Value = 84
$ws1.Cells.Item(71, 5).Value = "LP1912"
$ws1.Cells.Item(72, 1).Value = "06:57:11"
$ws1.Cells.Item(72, 2).Value = "08:54"
$ws1.Cells.Item(72, 3).Value = "17_ROMERO"
$ws1.Cells.Item(72, 4).Value = 117
$ws1.Cells.Item(72, 5).Value = "LP1912"
$ws1.Cells.Item(73, 1).Value = "07:19:37"
$ws1.Cells.Item(73, 2).Value = "09:01"
$ws1.Cells.Item(73, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(73, 4).Value = 102
$ws1.Cells.Item(73, 5).Value = "LP1912"
$ws1.Cells.Item(74, 1).Value = "07:19:37"
$ws1.Cells.Item(74, 2).Value = "09:10"
$ws1.Cells.Item(74, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(74, 4).Value = 111
$ws1.Cells.Item(74, 5).Value = "LP1912"
$ws1.Cells.Item(75, 1).Value = "07:19:37"
$ws1.Cells.Item(75, 2).Value = "09:16"
$ws1.Cells.Item(75, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(75, 4).Value = 117
$ws1.Cells.Item(75, 5).Value = "LP1912"
# Sheet 2: LP1912-215 -- update header + add row 17
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 07:19:37"
$ws2.Range("A3").Value = "Total filas: 12"
$ws2.Cells.Item(17, 1).Value = "07:19:37"
$ws2.Cells.Item(17, 2).Value = "09:01"
$ws2.Cells.Item(17, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(17, 4).Value = 102
$ws2.Cells.Item(17, 5).Value = "LP1912"
# Sheet 3: 6203-6173 -- update header + rows 14-20
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 07:19:37"
$ws3.Range("A3").Value = "Total filas: 15"
$ws3.Cells.Item(14, 1).Value = "07:19:37"
$ws3.Cells.Item(14, 2).Value = "07:38"
$ws3.Cells.Item(14, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(14, 4).Value = 19
$ws3.Cells.Item(14, 5).Value = "L6173"
$ws3.Cells.Item(15, 1).Value = "06:14:19"
$ws3.Cells.Item(15, 2).Value = "08:07"
$ws3.Cells.Item(15, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(15, 4).Value = 113
$ws3.Cells.Item(15, 5).Value = "L6203"
$ws3.Cells.Item(16, 1).Value = "07:19:37"
$ws3.Cells.Item(16, 2).Value = "08:08"
$ws3.Cells.Item(16, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(16, 4).Value = 49
$ws3.Cells.Item(16, 5).Value = "L6203"
$ws3.Cells.Item(17, 1).Value = "06:57:11"
$ws3.Cells.Item(17, 2).Value = "08:10"
$ws3.Cells.Item(17, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(17, 4).Value = 73
$ws3.Cells.Item(17, 5).Value = "L6203"
$ws3.Cells.Item(18, 1).Value = "07:19:37"
$ws3.Cells.Item(18, 2).Value = "08:35"
$ws3.Cells.Item(18, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18, 4).Value = 76
$ws3.Cells.Item(18, 5).Value = "L6173"
$ws3.Cells.Item(19, 1).Value = "06:57:11"
$ws3.Cells.Item(19, 2).Value = "08:38"
$ws3.Cells.Item(19, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(19, 4).Value = 101
$ws3.Cells.Item(19, 5).Value = "L6173"
$ws3.Cells.Item(20, 1).Value = "07:19:37"
$ws3.Cells.Item(20, 2).Value = "09:08"
$ws3.Cells.Item(20, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(20, 4).Value = 109
$ws3.Cells.Item(20, 5).Value = "L6203"
